$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.030015666666667
$ws.Cells.Item(2, 8).Value = 3.090047
$ws.Cells.Item(2, 9).Value = 0.001946685538032228
$ws.Cells.Item(2, 10).Value = 0.001946685538032228
$ws.Cells.Item(2, 13).Value = 3.795192333333334
$ws.Cells.Item(2, 14).Value = 11.385577
$ws.Cells.Item(2, 15).Value = 0.01044213755712683
$ws.Cells.Item(2, 16).Value = 0.01044213755712683
$ws.Cells.Item(2, 17).Value = 3.909107561346556
$ws.Cells.Item(2, 18).Value = 35.181968052119
$ws.Cells.Item(2, 19).Value = 0.00002032755816860199
$ws.Cells.Item(2, 20).Value = 0.00002032755816860199
$ws.Cells.Item(3, 7).Value = 1.030015666666667
$ws.Cells.Item(3, 8).Value = 3.090047
$ws.Cells.Item(3, 9).Value = 0.001946685538032228
$ws.Cells.Item(3, 10).Value = 0.001946685538032228
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.6696287328350964
$ws.Cells.Item(3, 16).Value = 0.6696287328350964
$ws.Cells.Item(3, 17).Value = 250.6815035235791
$ws.Cells.Item(3, 18).Value = 2256.133531712213
$ws.Cells.Item(3, 19).Value = 0.001303556570060929
$ws.Cells.Item(3, 20).Value = 0.001303556570060929
$ws.Cells.Item(4, 7).Value = 1.030015666666667
$ws.Cells.Item(4, 8).Value = 3.090047
$ws.Cells.Item(4, 9).Value = 0.001946685538032228
$ws.Cells.Item(4, 10).Value = 0.001946685538032228
$ws.Cells.Item(4, 13).Value = 29.801371
$ws.Cells.Item(4, 14).Value = 89.404113
$ws.Cells.Item(4, 15).Value = 0.08199584844219236
$ws.Cells.Item(4, 16).Value = 0.08199584844219235
$ws.Cells.Item(4, 17).Value = 30.69587901814567
$ws.Cells.Item(4, 18).Value = 276.262911163311
$ws.Cells.Item(4, 19).Value = 0.0001596201323410983
$ws.Cells.Item(4, 20).Value = 0.0001596201323410983
$ws.Cells.Item(5, 7).Value = 1.030015666666667
$ws.Cells.Item(5, 8).Value = 3.090047
$ws.Cells.Item(5, 9).Value = 0.001946685538032228
$ws.Cells.Item(5, 10).Value = 0.001946685538032228
$ws.Cells.Item(5, 13).Value = 86.47679266666667
$ws.Cells.Item(5, 14).Value = 259.430378
$ws.Cells.Item(5, 15).Value = 0.2379332811655844
$ws.Cells.Item(5, 16).Value = 0.2379332811655844
$ws.Cells.Item(5, 17).Value = 89.07245124975178
$ws.Cells.Item(5, 18).Value = 801.6520612477661
$ws.Cells.Item(5, 19).Value = 0.0004631812774615991
$ws.Cells.Item(5, 20).Value = 0.0004631812774615991
$ws.Cells.Item(6, 9).Value = 0.0008749538014921605
$ws.Cells.Item(6, 10).Value = 0.0008749538014921605
$ws.Cells.Item(6, 13).Value = 3.795192333333334
$ws.Cells.Item(6, 14).Value = 11.385577
$ws.Cells.Item(6, 15).Value = 0.01044213755712683
$ws.Cells.Item(6, 16).Value = 0.01044213755712683
$ws.Cells.Item(6, 17).Value = 1.756980495524334
$ws.Cells.Item(6, 18).Value = 15.812824459719
$ws.Cells.Item(6, 19).Value = 0.000009136387951312185
$ws.Cells.Item(6, 20).Value = 0.000009136387951312185
$ws.Cells.Item(7, 9).Value = 0.0008749538014921605
$ws.Cells.Item(7, 10).Value = 0.0008749538014921605
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.6696287328350964
$ws.Cells.Item(7, 16).Value = 0.6696287328350964
$ws.Cells.Item(7, 19).Value = 0.000585894205382446
$ws.Cells.Item(7, 20).Value = 0.000585894205382446
$ws.Cells.Item(8, 9).Value = 0.0008749538014921605
$ws.Cells.Item(8, 10).Value = 0.0008749538014921605
$ws.Cells.Item(8, 13).Value = 29.801371
$ws.Cells.Item(8, 14).Value = 89.404113
$ws.Cells.Item(8, 15).Value = 0.08199584844219236
$ws.Cells.Item(8, 16).Value = 0.08199584844219235
$ws.Cells.Item(8, 17).Value = 13.796514903079
$ws.Cells.Item(8, 18).Value = 124.168634127711
$ws.Cells.Item(8, 19).Value = 0.00007174257930107126
$ws.Cells.Item(8, 20).Value = 0.00007174257930107124
$ws.Cells.Item(9, 9).Value = 0.0008749538014921605
$ws.Cells.Item(9, 10).Value = 0.0008749538014921605
$ws.Cells.Item(9, 13).Value = 86.47679266666667
$ws.Cells.Item(9, 14).Value = 259.430378
$ws.Cells.Item(9, 15).Value = 0.2379332811655844
$ws.Cells.Item(9, 16).Value = 0.2379332811655844
$ws.Cells.Item(9, 17).Value = 40.03434468824067
$ws.Cells.Item(9, 18).Value = 360.309102194166
$ws.Cells.Item(9, 19).Value = 0.0002081806288573311
$ws.Cells.Item(9, 20).Value = 0.0002081806288573311
$ws.Cells.Item(10, 7).Value = 0.1159013333333333
$ws.Cells.Item(10, 8).Value = 0.347704
$ws.Cells.Item(10, 9).Value = 0.0002190485608522971
$ws.Cells.Item(10, 10).Value = 0.0002190485608522971
$ws.Cells.Item(10, 13).Value = 3.795192333333334
$ws.Cells.Item(10, 14).Value = 11.385577
$ws.Cells.Item(10, 15).Value = 0.01044213755712683
$ws.Cells.Item(10, 16).Value = 0.01044213755712683
$ws.Cells.Item(10, 17).Value = 0.4398678516897779
$ws.Cells.Item(10, 18).Value = 3.958810665208001
$ws.Cells.Item(10, 19).Value = 0.000002287335204110354
$ws.Cells.Item(10, 20).Value = 0.000002287335204110354
$ws.Cells.Item(11, 7).Value = 0.1159013333333333
$ws.Cells.Item(11, 8).Value = 0.347704
$ws.Cells.Item(11, 9).Value = 0.0002190485608522971
$ws.Cells.Item(11, 10).Value = 0.0002190485608522971
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.6696287328350964
$ws.Cells.Item(11, 16).Value = 0.6696287328350964
$ws.Cells.Item(11, 17).Value = 28.20764910733156
$ws.Cells.Item(11, 18).Value = 253.868841965984
$ws.Cells.Item(11, 19).Value = 0.0001466812102328752
$ws.Cells.Item(11, 20).Value = 0.0001466812102328752
$ws.Cells.Item(12, 7).Value = 0.1159013333333333
$ws.Cells.Item(12, 8).Value = 0.347704
$ws.Cells.Item(12, 9).Value = 0.0002190485608522971
$ws.Cells.Item(12, 10).Value = 0.0002190485608522971
$ws.Cells.Item(12, 13).Value = 29.801371
$ws.Cells.Item(12, 14).Value = 89.404113
$ws.Cells.Item(12, 15).Value = 0.08199584844219236
$ws.Cells.Item(12, 16).Value = 0.08199584844219235
$ws.Cells.Item(12, 17).Value = 3.454018634061334
$ws.Cells.Item(12, 18).Value = 31.086167706552
$ws.Cells.Item(12, 19).Value = 0.0000179610725971253
$ws.Cells.Item(12, 20).Value = 0.0000179610725971253
$ws.Cells.Item(13, 7).Value = 0.1159013333333333
$ws.Cells.Item(13, 8).Value = 0.347704
$ws.Cells.Item(13, 9).Value = 0.0002190485608522971
$ws.Cells.Item(13, 10).Value = 0.0002190485608522971
$ws.Cells.Item(13, 13).Value = 86.47679266666667
$ws.Cells.Item(13, 14).Value = 259.430378
$ws.Cells.Item(13, 15).Value = 0.2379332811655844
$ws.Cells.Item(13, 16).Value = 0.2379332811655844
$ws.Cells.Item(13, 17).Value = 10.02277557245689
$ws.Cells.Item(13, 18).Value = 90.20498015211201
$ws.Cells.Item(13, 19).Value = 0.00005211894281818622
$ws.Cells.Item(13, 20).Value = 0.00005211894281818622
$ws.Cells.Item(14, 7).Value = 527.5036416666667
$ws.Cells.Item(14, 8).Value = 1582.510925
$ws.Cells.Item(14, 9).Value = 0.9969593120996233
$ws.Cells.Item(14, 10).Value = 0.9969593120996233
$ws.Cells.Item(14, 13).Value = 3.795192333333334
$ws.Cells.Item(14, 14).Value = 11.385577
$ws.Cells.Item(14, 15).Value = 0.01044213755712683
$ws.Cells.Item(14, 16).Value = 0.01044213755712683
$ws.Cells.Item(14, 17).Value = 2001.977776658748
$ws.Cells.Item(14, 18).Value = 18017.79998992873
$ws.Cells.Item(14, 19).Value = 0.01041038627580281
$ws.Cells.Item(14, 20).Value = 0.01041038627580281
$ws.Cells.Item(15, 7).Value = 527.5036416666667
$ws.Cells.Item(15, 8).Value = 1582.510925
$ws.Cells.Item(15, 9).Value = 0.9969593120996233
$ws.Cells.Item(15, 10).Value = 0.9969593120996233
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.6696287328350964
$ws.Cells.Item(15, 16).Value = 0.6696287328350964
$ws.Cells.Item(15, 17).Value = 128381.9365923852
$ws.Cells.Item(15, 18).Value = 1155437.429331467
$ws.Cells.Item(15, 19).Value = 0.6675926008494201
$ws.Cells.Item(15, 20).Value = 0.6675926008494201
$ws.Cells.Item(16, 7).Value = 527.5036416666667
$ws.Cells.Item(16, 8).Value = 1582.510925
$ws.Cells.Item(16, 9).Value = 0.9969593120996233
$ws.Cells.Item(16, 10).Value = 0.9969593120996233
$ws.Cells.Item(16, 13).Value = 29.801371
$ws.Cells.Item(16, 14).Value = 89.404113
$ws.Cells.Item(16, 15).Value = 0.08199584844219236
$ws.Cells.Item(16, 16).Value = 0.08199584844219235
$ws.Cells.Item(16, 17).Value = 15720.33172915939
$ws.Cells.Item(16, 18).Value = 141482.9855624345
$ws.Cells.Item(16, 19).Value = 0.08174652465795307
$ws.Cells.Item(16, 20).Value = 0.08174652465795305
$ws.Cells.Item(17, 7).Value = 527.5036416666667
$ws.Cells.Item(17, 8).Value = 1582.510925
$ws.Cells.Item(17, 9).Value = 0.9969593120996233
$ws.Cells.Item(17, 10).Value = 0.9969593120996233
$ws.Cells.Item(17, 13).Value = 86.47679266666667
$ws.Cells.Item(17, 14).Value = 259.430378
$ws.Cells.Item(17, 15).Value = 0.2379332811655844
$ws.Cells.Item(17, 16).Value = 0.2379332811655844
$ws.Cells.Item(17, 17).Value = 45616.82305131996
$ws.Cells.Item(17, 18).Value = 410551.4074618797
$ws.Cells.Item(17, 19).Value = 0.2372098003164473
$ws.Cells.Item(17, 20).Value = 0.2372098003164473
